$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 112, pushing existing rows 112:223 down to 113:224
$ws.Rows.Item(112).Insert()

# Populate the newly inserted row 112 with the new data record
$ws.Range("A112").Value = 3
$ws.Range("B112").Value = "Femacal de La Calera"
$ws.Range("C112").Value = "Coquimbo"
$ws.Range("D112").Value = 44587
$ws.Range("E112").Value = 5
$ws.Range("F112").Value = 100112001
$ws.Range("G112").Value = "Berenjena"
$ws.Range("H112").Value = "Sin especificar"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 65
$ws.Range("K112").Value = 9500
$ws.Range("L112").Value = 10000
$ws.Range("M112").Value = 9731
$ws.Range("N112").Value = "$/caja 60 unidades"
$ws.Range("O112").Value = "Región de Arica y Parinacota"
$ws.Range("P112").Value = 162
$ws.Range("Q112").Value = 60
$ws.Range("R112").Value = "Hortaliza"
